$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused 4th row (table shrinks from 4 data-ish rows to a
# 1 header + 2 data row layout).
$ws.Rows.Item(4).Delete()

# Rewrite the remaining three rows with the new header + new trait text.
$ws.Range("A1").Value = "trait"
$ws.Range("B1").Value = "formula"
$ws.Range("A2").Value = "IgGI_first_trait"
$ws.Range("B2").Value = "0.5 * (IgGI1H4N4F1 + IgGI1H5N4F1)"
$ws.Range("A3").Value = "second_trait"
$ws.Range("B3").Value = "IgGI1H4N4F1S1 / (IgGI1H4N4F1 + IgGI1H4N5F1S1)"

# Header row becomes bold, matching the default formatting used for the
# rest of the generated traits tables.
$ws.Range("A1:B1").Font.Bold = $true

# Clear the old manual selection/active cell so the sheet opens at the
# default view.
[void]$ws.Range("A1").Select()

# Widen the columns so the longer trait/formula text is readable.
$ws.Columns.Item(1).ColumnWidth = 14.8932291667
$ws.Columns.Item(2).ColumnWidth = 41.6197916667

# Match default page setup (portrait, paper size 9 = A4) used elsewhere.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
